$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-looking decimal numbers as text
# (e.g. "1.001", "307.40", "0.000008546"). Excel's Range.Value setter
# auto-converts numeric-looking strings to real numbers, which would
# lose the original text formatting (trailing zeros, thousand-dot
# grouping, etc). Force those specific cells to stay text first.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.167.89'
$ws.Range("E2").Value = '  -2.31%  '
$ws.Range("D3").Value = '1.871.05'
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '307.40'
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '0.5110'
$ws.Range("E7").Value = '  +1.45%  '
$ws.Range("D8").Value = '0.3762'
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").Value = '0.07169'
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").Value = '0.8901'
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("D11").Value = '20.72'
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").Value = '0.07599'
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").Value = '1.868.51'
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("D14").Value = '5.336'
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("D15").Value = '89.37'
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '0.000008546'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").Value = '14.15'
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '27.206.98'
$ws.Range("E20").Value = '  -2.27%  '
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").Value = '2.098.94'
$ws.Range("E22").Value = '  -2.25%  '
$ws.Range("D23").Value = '10.63'
$ws.Range("E23").Value = '  -1.61%  '
$ws.Range("D24").Value = '6.491'
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("D25").Value = '150.84'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").Value = '1.839'
$ws.Range("E26").Value = '  -1.88%  '
$ws.Range("D27").Value = '18.03'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("D28").Value = '2.120'
$ws.Range("E28").Value = '  -4.30%  '
$ws.Range("D29").Value = '112.89'
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("D30").Value = '4.757'
$ws.Range("D31").Value = '4.720'
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("D32").Value = '0.09004'
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("D34").Value = '3.096'
$ws.Range("E34").Value = '  -3.63%  '
$ws.Range("D35").Value = '0.7538'
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").Value = '1.168'
$ws.Range("E36").Value = '  -3.97%  '
$ws.Range("D37").Value = '0.02032'
$ws.Range("E37").Value = '  -1.39%  '
$ws.Range("D38").Value = '2.532'
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("D39").Value = '3.028'
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("D40").Value = '1.075'
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("D41").Value = '0.5355'
$ws.Range("E41").Value = '  -2.73%  '
$ws.Range("D42").Value = '6.643'
$ws.Range("E42").Value = '  -3.02%  '
$ws.Range("D43").Value = '113.86'
$ws.Range("D44").Value = '8.555'
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("D45").Value = '0.1481'
$ws.Range("E45").Value = '  -1.77%  '
$ws.Range("D46").Value = '0.4668'
$ws.Range("E46").Value = '  -2.57%  '
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").Value = '10.12'
$ws.Range("E48").Value = '  -4.20%  '
$ws.Range("D49").Value = '1.571'
$ws.Range("E49").Value = '  -3.07%  '
$ws.Range("D50").Value = '65.14'
$ws.Range("E50").Value = '  -3.02%  '
$ws.Range("D51").Value = '36.70'
$ws.Range("E51").Value = '  -0.43%  '
